$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 138, shifting existing rows 138:248 down to 139:249
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new record's data
$ws.Cells.Item(138, 1).Value = 11
$ws.Cells.Item(138, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(138, 3).Value = "Bíobío"
$ws.Cells.Item(138, 4).Value = 44957
$ws.Cells.Item(138, 5).Value = 8
$ws.Cells.Item(138, 6).Value = 100112003
$ws.Cells.Item(138, 7).Value = "Ajo"
$ws.Cells.Item(138, 8).Value = "Chino"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 180
$ws.Cells.Item(138, 11).Value = 15000
$ws.Cells.Item(138, 12).Value = 16000
$ws.Cells.Item(138, 13).Value = 15556
$ws.Cells.Item(138, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(138, 15).Value = "China"
$ws.Cells.Item(138, 16).Value = 1556
$ws.Cells.Item(138, 17).Value = 10
$ws.Cells.Item(138, 18).Value = "Hortaliza"
